# #5: cash & deposit done
# Populate the 存款 (deposit) sheet with bank/deposit_type/currency columns
# plus the common trailing metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- header row (row 1) ------------------------------------------------
# B1/C1/D1 (bank/deposit_type/currency) already hold the right labels.
# E1/F1 need to switch from literal sample values to the owner/total
# column headers, and G1:M1 are brand new header cells.
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "total"
$ws.Cells.Item(1,7).Value  = "property_category"
$ws.Cells.Item(1,8).Value  = "category"
$ws.Cells.Item(1,9).Value  = "date"
$ws.Cells.Item(1,10).Value = "legislator_name"
$ws.Cells.Item(1,11).Value = "legislator_id"
$ws.Cells.Item(1,12).Value = "source_file"
$ws.Cells.Item(1,13).Value = "index"

# match the bold/centered/bordered header style used by B1:F1
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null

# --- data rows (rows 2-8) ----------------------------------------------
$bank = "國泰世華商業銀行仁愛分行"
$bankFubon = "台北富邦商業銀行玉成分行"
$bankCoop = "合作金庫商業銀行建國分行"

$fixed = "定期存款"
$demand = "活期存款"
$check = "支票存款"

$currency = "新臺幣"

$owner1 = "吳宜臻"
$owner2 = "楊〇翰"
$owner3 = "楊◦霓"

$rows = @(
    @{ Row=2; Index=45; Bank=$bank;       Type=$fixed;  Owner=$owner1; Total=200000 },
    @{ Row=3; Index=46; Bank=$bankFubon;  Type=$demand; Owner=$owner1; Total=499424 },
    @{ Row=4; Index=47; Bank=$bank;       Type=$demand; Owner=$owner1; Total=453908 },
    @{ Row=5; Index=48; Bank=$bank;       Type=$check;  Owner=$owner1; Total=144522 },
    @{ Row=6; Index=49; Bank=$bank;       Type=$demand; Owner=$owner2; Total=809749 },
    @{ Row=7; Index=50; Bank=$bankCoop;   Type=$demand; Owner=$owner2; Total=366838 },
    @{ Row=8; Index=51; Bank=$bank;       Type=$demand; Owner=$owner3; Total=661096 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row,2).Value = $r.Bank
    $ws.Cells.Item($row,3).Value = $r.Type
    $ws.Cells.Item($row,4).Value = $currency
    $ws.Cells.Item($row,5).Value = $r.Owner
    $ws.Cells.Item($row,6).Value = $r.Total

    $ws.Cells.Item($row,7).Value  = "deposit"
    $ws.Cells.Item($row,8).Value  = "normal"

    # force the ISO date to stay a plain text string instead of being
    # auto-converted to a date serial number
    $ws.Cells.Item($row,9).NumberFormat = "@"
    $ws.Cells.Item($row,9).Value  = "2012-02-01"

    $ws.Cells.Item($row,10).Value = "吳宜臻"
    $ws.Cells.Item($row,11).Value = 1735
    $ws.Cells.Item($row,12).Value = "tmp2691"
    $ws.Cells.Item($row,13).Value = $r.Index
}

$excel.CutCopyMode = 0
